$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "321.70"
Set-TextValue "E2" "8.69%"
Set-TextValue "G2" "2"
Set-TextValue "D3" "46.65"
Set-TextValue "E3" "11.02%"
Set-TextValue "G3" "2"
Set-TextValue "D4" "5.264"
Set-TextValue "E4" "4.99%"
Set-TextValue "G4" "2"
Set-TextValue "D5" "0.08106"
Set-TextValue "E5" "7.71%"
Set-TextValue "G5" "2"
Set-TextValue "D6" "4.551"
Set-TextValue "E6" "3.68%"
Set-TextValue "G6" "2"
Set-TextValue "D7" "1.691"
Set-TextValue "E7" "7.48%"
Set-TextValue "G7" "2"
Set-TextValue "D8" "1.093"
Set-TextValue "E8" "18.07%"
Set-TextValue "G8" "2"
Set-TextValue "D9" "0.1310"
Set-TextValue "E9" "10.52%"
Set-TextValue "G9" "2"
Set-TextValue "D10" "0.1962"
Set-TextValue "E10" "6.83%"
Set-TextValue "G10" "2"
Set-TextValue "D11" "0.09624"
Set-TextValue "E11" "6.82%"
Set-TextValue "G11" "2"
Set-TextValue "D12" "0.04376"
Set-TextValue "E12" "4.56%"
Set-TextValue "G12" "2"
Set-TextValue "D13" "0.1047"
Set-TextValue "E13" "-0.18%"
Set-TextValue "G13" "2"
Set-TextValue "D14" "0.001323"
Set-TextValue "E14" "3.23%"
Set-TextValue "G14" "2"
Set-TextValue "D15" "0.005868"
Set-TextValue "E15" "0.32%"
Set-TextValue "G15" "2"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.418"
Set-TextValue "E16" "1.71%"
Set-TextValue "G16" "2"
Set-TextValue "B17" "BTSEToken"
Set-TextValue "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D17" "2.437"
Set-TextValue "E17" "1.51%"
Set-TextValue "G17" "2"
Set-TextValue "B18" "BitpandaEcosystemToken"
Set-TextValue "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D18" "0.3390"
Set-TextValue "E18" "2.22%"
Set-TextValue "G18" "2"
Set-TextValue "B19" "MCDex"
Set-TextValue "C19" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D19" "8.170"
Set-TextValue "E19" "2.01%"
Set-TextValue "G19" "2"
Set-TextValue "B20" "ProBitToken"
Set-TextValue "C20" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D20" "0.1398"
Set-TextValue "E20" "1.07%"
Set-TextValue "G20" "2"
Set-TextValue "B21" "ZBToken"
Set-TextValue "C21" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D21" "0.3073"
Set-TextValue "E21" "2.27%"
Set-TextValue "G21" "2"
Set-TextValue "B22" "CoinExToken"
Set-TextValue "C22" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D22" "0.04327"
Set-TextValue "E22" "4.63%"
Set-TextValue "G22" "2"
Set-TextValue "B23" "BitKan"
Set-TextValue "C23" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D23" "0.001299"
Set-TextValue "E23" "2.57%"
Set-TextValue "G23" "2"
Set-TextValue "B24" "HotbitToken"
Set-TextValue "C24" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D24" "0.004261"
Set-TextValue "E24" "8.32%"
Set-TextValue "G24" "2"
Set-TextValue "B25" "NitroEx"
Set-TextValue "C25" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D25" "0.0001346"
Set-TextValue "E25" "9.27%"
Set-TextValue "G25" "2"
Set-TextValue "B26" "UpBots"
Set-TextValue "C26" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D26" "0.0003711"
Set-TextValue "E26" "-0.34%"
Set-TextValue "G26" "2"
Set-TextValue "G27" "2"
Set-TextValue "G28" "2"
Set-TextValue "G29" "2"
Set-TextValue "G30" "2"
Set-TextValue "G31" "2"
Set-TextValue "G32" "2"
Set-TextValue "G33" "2"
Set-TextValue "G34" "2"
Set-TextValue "G35" "2"
Set-TextValue "G36" "2"
Set-TextValue "G37" "2"
Set-TextValue "D38" "0.02751"
Set-TextValue "E38" "14.59%"
Set-TextValue "G38" "2"
Set-TextValue "D39" "0.05534"
Set-TextValue "E39" "6.35%"
Set-TextValue "G39" "2"
Set-TextValue "D40" "0.005872"
Set-TextValue "E40" "-6.98%"
Set-TextValue "G40" "2"
Set-TextValue "D41" "0.007762"
Set-TextValue "E41" "-0.52%"
Set-TextValue "G41" "2"
Set-TextValue "D42" "0.1447"
Set-TextValue "E42" "9.19%"
Set-TextValue "G42" "2"
Set-TextValue "D43" "0.007653"
Set-TextValue "E43" "3.16%"
Set-TextValue "G43" "2"
Set-TextValue "D44" "0.008924"
Set-TextValue "E44" "23.43%"
Set-TextValue "G44" "2"
Set-TextValue "D45" "0.3218"
Set-TextValue "E45" "0.39%"
Set-TextValue "G45" "2"
Set-TextValue "D46" "0.00006848"
Set-TextValue "E46" "6.68%"
Set-TextValue "G46" "2"
Set-TextValue "E47" "-0.49%"
Set-TextValue "G47" "2"
Set-TextValue "D48" "0.05330"
Set-TextValue "E48" "48.96%"
Set-TextValue "G48" "2"
Set-TextValue "D49" "0.003989"
Set-TextValue "E49" "-5.22%"
Set-TextValue "G49" "2"
Set-TextValue "E50" "-0.49%"
Set-TextValue "G50" "2"
Set-TextValue "E51" "-0.49%"
Set-TextValue "G51" "2"
